# Update activity tracking data for three rows in the Banco_Dashboard sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# Row 5: "Abrir BV laterais parte cilindrica SDA" - now fully completed
$ws.Range("D5").Value = 10
$ws.Range("F2").Copy()
$ws.Range("G5").PasteSpecial($xlPasteFormats)
$ws.Range("G5").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("H5").Value = 100

# Row 40: "Abrir 12 BV inferiores hoppers" - now in progress
$ws.Range("D40").Value = 5
$ws.Range("F2").Copy()
$ws.Range("F40").PasteSpecial($xlPasteFormats)
$ws.Range("F40").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("H40").Value = 40

# Row 61: "Abrir 01 (uma) boca de visita que dá acesso ao duto de saída dos Idfan's" - now fully completed
$ws.Range("D61").Value = 1
$ws.Range("F2").Copy()
$ws.Range("F61").PasteSpecial($xlPasteFormats)
$ws.Range("F61").Value = (Get-Date -Year 2025 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("F2").Copy()
$ws.Range("G61").PasteSpecial($xlPasteFormats)
$ws.Range("G61").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("H61").Value = 100

$excel.CutCopyMode = $false

# Update the active selection/view to match the saved state
$ws.Range("H62").Select()

$wb.Save()
